$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.212.05"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.479.83"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.78"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.54"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "2.479.16"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.56"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "2.929.77"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "63.192.26"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "2.478.71"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.25"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  +8.83%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.19"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  +15.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "658.76"
$ws.Range("E27").Value = "  +5.69%  "
$ws.Range("D28").Value = "0.0₃0992"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +281.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +5.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.132"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.46"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.31"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "0.0₆0318"
$ws.Range("E45").Value = "  -47.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.22"
$ws.Range("E46").Value = "  +7.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.24"
$ws.Range("E47").Value = "  +3.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.61"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.43"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("E51").Value = "  +0.17%  "
